# Auto-generated: apply value updates per the commit diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 550.875
$ws.Range("J17").Value = 562.2258
$ws.Range("L17").Value = 1686.6774
$ws.Range("N17").Value = -2022.6774
$ws.Range("H113").Value = 3110.7144
$ws.Range("I113").Value = 3075
$ws.Range("K113").Value = 3075
$ws.Range("M113").Value = 179
$ws.Range("H116").Value = 6582.263
$ws.Range("I116").Value = 5874.9165
$ws.Range("K116").Value = 5874.9165
$ws.Range("M116").Value = -2432.9165
$ws.Range("H132").Value = 2043.24
$ws.Range("I132").Value = 2134.087
$ws.Range("K132").Value = 6402.261
$ws.Range("M132").Value = -3872.261
$ws.Range("H135").Value = 2176.1052
$ws.Range("I135").Value = 1959.1428
$ws.Range("K135").Value = 17632.2852
$ws.Range("M135").Value = -15097.2852
$ws.Range("H141").Value = 7848.727
$ws.Range("I141").Value = 7042.125
$ws.Range("K141").Value = 21126.375
$ws.Range("M141").Value = -15946.375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5954
$ws.Range("I32").Value = 2749.8774
$ws.Range("K32").Value = 2749.8774
$ws.Range("M32").Value = -2462.8774
$ws.Range("H61").Value = 1301.0385
$ws.Range("I61").Value = 1014.86365
$ws.Range("K61").Value = 1014.86365
$ws.Range("M61").Value = -802.86365
$ws.Range("H74").Value = 2046.421
$ws.Range("I74").Value = 1552.3572
$ws.Range("J74").Value = 3429.8
$ws.Range("K74").Value = 1552.3572
$ws.Range("L74").Value = 3429.8
$ws.Range("M74").Value = -678.3571999999999
$ws.Range("N74").Value = -5177.8
$ws.Range("H77").Value = 2046.421
$ws.Range("I77").Value = 1552.3572
$ws.Range("J77").Value = 3429.8
$ws.Range("K77").Value = 7761.786
$ws.Range("L77").Value = 17149
$ws.Range("M77").Value = -3393.786
$ws.Range("N77").Value = -25885
$ws.Range("H132").Value = 1338.5358
$ws.Range("I132").Value = 1222.24
$ws.Range("K132").Value = 3666.72
$ws.Range("M132").Value = -1136.72
$ws.Range("H136").Value = 1301.0385
$ws.Range("I136").Value = 1014.86365
$ws.Range("K136").Value = 3044.59095
$ws.Range("M136").Value = -494.5909499999998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1050.5
$ws.Range("I5").Value = 850
$ws.Range("J5").Value = 1150.75
$ws.Range("K5").Value = 850
$ws.Range("L5").Value = 1150.75
$ws.Range("M5").Value = -737
$ws.Range("N5").Value = -1376.75
$ws.Range("H20").Value = 720.37036
$ws.Range("I20").Value = 684.3889
$ws.Range("J20").Value = 792.3333
$ws.Range("K20").Value = 684.3889
$ws.Range("L20").Value = 792.3333
$ws.Range("M20").Value = -437.3889
$ws.Range("N20").Value = -1286.3333
$ws.Range("H22").Value = 3119819
$ws.Range("I22").Value = 3898929
$ws.Range("K22").Value = 3898929
$ws.Range("M22").Value = -3898756
$ws.Range("H86").Value = 12167.167
$ws.Range("I86").Value = 8249
$ws.Range("K86").Value = 8249
$ws.Range("M86").Value = -7126
$ws.Range("H89").Value = 12167.167
$ws.Range("I89").Value = 8249
$ws.Range("K89").Value = 41245
$ws.Range("M89").Value = -35629
$ws.Range("H134").Value = 9449.615
$ws.Range("I134").Value = 9418.125
$ws.Range("J134").Value = 9500
$ws.Range("K134").Value = 28254.375
$ws.Range("L134").Value = 28500
$ws.Range("M134").Value = -25719.375
$ws.Range("N134").Value = -33570

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 25000
$ws.Range("J57").Value = 25000
$ws.Range("L57").Value = 25000
$ws.Range("N57").Value = -26120
$ws.Range("H99").Value = 6175225
$ws.Range("I99").Value = 9261283
$ws.Range("J99").Value = 3108.1667
$ws.Range("K99").Value = 9261283
$ws.Range("L99").Value = 3108.1667
$ws.Range("M99").Value = -9259785
$ws.Range("N99").Value = -6104.1667
$ws.Range("H126").Value = 6175225
$ws.Range("I126").Value = 9261283
$ws.Range("J126").Value = 3108.1667
$ws.Range("K126").Value = 27783849
$ws.Range("L126").Value = 9324.500100000001
$ws.Range("M126").Value = -27781379
$ws.Range("N126").Value = -14264.5001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 201.25
$ws.Range("I33").Value = 253.14285
$ws.Range("J33").Value = 128.6
$ws.Range("K33").Value = 1518.8571
$ws.Range("L33").Value = 771.5999999999999
$ws.Range("M33").Value = -1235.8571
$ws.Range("N33").Value = -1337.6
$ws.Range("H123").Value = 3899.6667
$ws.Range("J123").Value = 9999
$ws.Range("L123").Value = 29997
$ws.Range("N123").Value = -34897

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6823
$ws.Range("J70").Value = 6739.8
$ws.Range("L70").Value = 6739.8
$ws.Range("N70").Value = -7279.8
$ws.Range("H73").Value = 6823
$ws.Range("J73").Value = 6739.8
$ws.Range("L73").Value = 6739.8
$ws.Range("N73").Value = -8611.799999999999
$ws.Range("H124").Value = 98000
$ws.Range("J124").Value = 98000
$ws.Range("L124").Value = 98000
$ws.Range("N124").Value = -107820
$ws.Range("H126").Value = 3695.9583
$ws.Range("I126").Value = 2151.4167
$ws.Range("K126").Value = 6454.250100000001
$ws.Range("M126").Value = -3984.250100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1369.3846
$ws.Range("I22").Value = 1150
$ws.Range("K22").Value = 1150
$ws.Range("M22").Value = -855
$ws.Range("H27").Value = 1369.3846
$ws.Range("I27").Value = 1150
$ws.Range("K27").Value = 1150
$ws.Range("M27").Value = -1043
$ws.Range("H46").Value = 3607.8696
$ws.Range("I46").Value = 2600
$ws.Range("J46").Value = 3887.8333
$ws.Range("K46").Value = 2600
$ws.Range("L46").Value = 3887.8333
$ws.Range("M46").Value = -2412
$ws.Range("N46").Value = -4263.8333
$ws.Range("H55").Value = 7450
$ws.Range("I55").Value = 6000
$ws.Range("J55").Value = 8175
$ws.Range("K55").Value = 6000
$ws.Range("L55").Value = 8175
$ws.Range("M55").Value = -5827
$ws.Range("N55").Value = -8521
$ws.Range("H61").Value = 1121.7142
$ws.Range("I61").Value = 1121.7142
$ws.Range("K61").Value = 1121.7142
$ws.Range("M61").Value = -919.7141999999999
$ws.Range("H113").Value = 1121.7142
$ws.Range("I113").Value = 1121.7142
$ws.Range("K113").Value = 1121.7142
$ws.Range("M113").Value = 1048.2858
$ws.Range("H136").Value = 2554.0688
$ws.Range("I136").Value = 3116
$ws.Range("J136").Value = 1862.4615
$ws.Range("K136").Value = 9348
$ws.Range("L136").Value = 5587.3845
$ws.Range("M136").Value = -6798
$ws.Range("N136").Value = -10687.3845

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4374.75
$ws.Range("J81").Value = 4374.75
$ws.Range("L81").Value = 8749.5
$ws.Range("N81").Value = -10871.5
$ws.Range("H84").Value = 4374.75
$ws.Range("J84").Value = 4374.75
$ws.Range("L84").Value = 43747.5
$ws.Range("N84").Value = -54355.5
$ws.Range("H132").Value = 1174.8077
$ws.Range("I132").Value = 1001.875
$ws.Range("K132").Value = 3005.625
$ws.Range("M132").Value = -475.625
